$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(549, 4).Value = 45131
$ws.Cells.Item(549, 9).Value = "Primera"
$ws.Cells.Item(549, 10).Value = 5000
$ws.Cells.Item(549, 11).Value = 550
$ws.Cells.Item(549, 12).Value = 550
$ws.Cells.Item(549, 13).Value = 550
$ws.Cells.Item(549, 16).Value = 110

$ws.Cells.Item(550, 4).Value = 44188
$ws.Cells.Item(550, 9).Value = "Primera"
$ws.Cells.Item(550, 10).Value = 3000
$ws.Cells.Item(550, 11).Value = 600
$ws.Cells.Item(550, 12).Value = 600
$ws.Cells.Item(550, 13).Value = 600
$ws.Cells.Item(550, 16).Value = 120

$ws.Cells.Item(551, 4).Value = 44376
$ws.Cells.Item(551, 9).Value = "Primera"
$ws.Cells.Item(551, 10).Value = 2000
$ws.Cells.Item(551, 11).Value = 600
$ws.Cells.Item(551, 12).Value = 600
$ws.Cells.Item(551, 13).Value = 600
$ws.Cells.Item(551, 16).Value = 120

$ws.Cells.Item(552, 4).Value = 44273
$ws.Cells.Item(552, 9).Value = "Primera"
$ws.Cells.Item(552, 10).Value = 3000
$ws.Cells.Item(552, 11).Value = 500
$ws.Cells.Item(552, 12).Value = 500
$ws.Cells.Item(552, 13).Value = 500
$ws.Cells.Item(552, 16).Value = 100

$ws.Cells.Item(553, 4).Value = 44215
$ws.Cells.Item(553, 9).Value = "Primera"
$ws.Cells.Item(553, 10).Value = 2000
$ws.Cells.Item(553, 11).Value = 600
$ws.Cells.Item(553, 12).Value = 600
$ws.Cells.Item(553, 13).Value = 600
$ws.Cells.Item(553, 16).Value = 120

$ws.Cells.Item(554, 4).Value = 45124
$ws.Cells.Item(554, 9).Value = "Primera"
$ws.Cells.Item(554, 10).Value = 5000
$ws.Cells.Item(554, 11).Value = 600
$ws.Cells.Item(554, 12).Value = 600
$ws.Cells.Item(554, 13).Value = 600
$ws.Cells.Item(554, 16).Value = 120

$ws.Cells.Item(555, 4).Value = 44902
$ws.Cells.Item(555, 9).Value = "Primera"
$ws.Cells.Item(555, 10).Value = 5000
$ws.Cells.Item(555, 11).Value = 700
$ws.Cells.Item(555, 12).Value = 700
$ws.Cells.Item(555, 13).Value = 700
$ws.Cells.Item(555, 16).Value = 140

$ws.Cells.Item(556, 4).Value = 44631
$ws.Cells.Item(556, 9).Value = "Segunda"
$ws.Cells.Item(556, 10).Value = 3000
$ws.Cells.Item(556, 11).Value = 600
$ws.Cells.Item(556, 12).Value = 600
$ws.Cells.Item(556, 13).Value = 600
$ws.Cells.Item(556, 16).Value = 120

$ws.Cells.Item(557, 4).Value = 44945
$ws.Cells.Item(557, 9).Value = "Primera"
$ws.Cells.Item(557, 10).Value = 3000
$ws.Cells.Item(557, 11).Value = 700
$ws.Cells.Item(557, 12).Value = 700
$ws.Cells.Item(557, 13).Value = 700
$ws.Cells.Item(557, 16).Value = 140

$ws.Cells.Item(558, 4).Value = 44554
$ws.Cells.Item(558, 9).Value = "Primera"
$ws.Cells.Item(558, 10).Value = 3000
$ws.Cells.Item(558, 11).Value = 500
$ws.Cells.Item(558, 12).Value = 500
$ws.Cells.Item(558, 13).Value = 500
$ws.Cells.Item(558, 16).Value = 100

$ws.Cells.Item(559, 4).Value = 45070
$ws.Cells.Item(559, 9).Value = "Primera"
$ws.Cells.Item(559, 10).Value = 3000
$ws.Cells.Item(559, 11).Value = 600
$ws.Cells.Item(559, 12).Value = 600
$ws.Cells.Item(559, 13).Value = 600
$ws.Cells.Item(559, 16).Value = 120

$ws.Cells.Item(560, 4).Value = 45070
$ws.Cells.Item(560, 9).Value = "Segunda"
$ws.Cells.Item(560, 10).Value = 2000
$ws.Cells.Item(560, 11).Value = 500
$ws.Cells.Item(560, 12).Value = 500
$ws.Cells.Item(560, 13).Value = 500
$ws.Cells.Item(560, 16).Value = 100

$ws.Cells.Item(561, 4).Value = 44790
$ws.Cells.Item(561, 9).Value = "Primera"
$ws.Cells.Item(561, 10).Value = 3000
$ws.Cells.Item(561, 11).Value = 750
$ws.Cells.Item(561, 12).Value = 750
$ws.Cells.Item(561, 13).Value = 750
$ws.Cells.Item(561, 16).Value = 150

$ws.Cells.Item(562, 4).Value = 44714
$ws.Cells.Item(562, 9).Value = "Primera"
$ws.Cells.Item(562, 10).Value = 3000
$ws.Cells.Item(562, 11).Value = 700
$ws.Cells.Item(562, 12).Value = 700
$ws.Cells.Item(562, 13).Value = 700
$ws.Cells.Item(562, 16).Value = 140

$ws.Cells.Item(563, 4).Value = 44272
$ws.Cells.Item(563, 9).Value = "Primera"
$ws.Cells.Item(563, 10).Value = 3000
$ws.Cells.Item(563, 11).Value = 500
$ws.Cells.Item(563, 12).Value = 500
$ws.Cells.Item(563, 13).Value = 500
$ws.Cells.Item(563, 16).Value = 100

$ws.Cells.Item(564, 4).Value = 45040
$ws.Cells.Item(564, 9).Value = "Primera"
$ws.Cells.Item(564, 10).Value = 3000
$ws.Cells.Item(564, 11).Value = 700
$ws.Cells.Item(564, 12).Value = 700
$ws.Cells.Item(564, 13).Value = 700
$ws.Cells.Item(564, 16).Value = 140

$ws.Cells.Item(565, 4).Value = 45040
$ws.Cells.Item(565, 9).Value = "Segunda"
$ws.Cells.Item(565, 10).Value = 2000
$ws.Cells.Item(565, 11).Value = 600
$ws.Cells.Item(565, 12).Value = 600
$ws.Cells.Item(565, 13).Value = 600
$ws.Cells.Item(565, 16).Value = 120

$ws.Cells.Item(566, 4).Value = 44826
$ws.Cells.Item(566, 9).Value = "Segunda"
$ws.Cells.Item(566, 10).Value = 3000
$ws.Cells.Item(566, 11).Value = 800
$ws.Cells.Item(566, 12).Value = 800
$ws.Cells.Item(566, 13).Value = 800
$ws.Cells.Item(566, 16).Value = 160

$ws.Cells.Item(567, 4).Value = 44692
$ws.Cells.Item(567, 9).Value = "Primera"
$ws.Cells.Item(567, 10).Value = 5000
$ws.Cells.Item(567, 11).Value = 600
$ws.Cells.Item(567, 12).Value = 600
$ws.Cells.Item(567, 13).Value = 600
$ws.Cells.Item(567, 16).Value = 120

$ws.Cells.Item(568, 4).Value = 44589
$ws.Cells.Item(568, 9).Value = "Primera"
$ws.Cells.Item(568, 10).Value = 4000
$ws.Cells.Item(568, 11).Value = 700
$ws.Cells.Item(568, 12).Value = 700
$ws.Cells.Item(568, 13).Value = 700
$ws.Cells.Item(568, 16).Value = 140

$ws.Cells.Item(569, 4).Value = 44771
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 3000
$ws.Cells.Item(569, 11).Value = 750
$ws.Cells.Item(569, 12).Value = 750
$ws.Cells.Item(569, 13).Value = 750
$ws.Cells.Item(569, 16).Value = 150

$ws.Cells.Item(570, 4).Value = 44565
$ws.Cells.Item(570, 9).Value = "Primera"
$ws.Cells.Item(570, 10).Value = 4000
$ws.Cells.Item(570, 11).Value = 500
$ws.Cells.Item(570, 12).Value = 500
$ws.Cells.Item(570, 13).Value = 500
$ws.Cells.Item(570, 16).Value = 100

$ws.Cells.Item(571, 4).Value = 44901
$ws.Cells.Item(571, 9).Value = "Primera"
$ws.Cells.Item(571, 10).Value = 5000
$ws.Cells.Item(571, 11).Value = 700
$ws.Cells.Item(571, 12).Value = 700
$ws.Cells.Item(571, 13).Value = 700
$ws.Cells.Item(571, 16).Value = 140

$ws.Cells.Item(572, 4).Value = 44322
$ws.Cells.Item(572, 9).Value = "Primera"
$ws.Cells.Item(572, 10).Value = 4000
$ws.Cells.Item(572, 11).Value = 500
$ws.Cells.Item(572, 12).Value = 500
$ws.Cells.Item(572, 13).Value = 500
$ws.Cells.Item(572, 16).Value = 100

$ws.Cells.Item(573, 4).Value = 44495
$ws.Cells.Item(573, 9).Value = "Primera"
$ws.Cells.Item(573, 10).Value = 4000
$ws.Cells.Item(573, 11).Value = 650
$ws.Cells.Item(573, 12).Value = 650
$ws.Cells.Item(573, 13).Value = 650
$ws.Cells.Item(573, 16).Value = 130

$ws.Cells.Item(574, 4).Value = 45111
$ws.Cells.Item(574, 9).Value = "Primera"
$ws.Cells.Item(574, 10).Value = 5000
$ws.Cells.Item(574, 11).Value = 600
$ws.Cells.Item(574, 12).Value = 600
$ws.Cells.Item(574, 13).Value = 600
$ws.Cells.Item(574, 16).Value = 120

$ws.Cells.Item(575, 4).Value = 44417
$ws.Cells.Item(575, 9).Value = "Primera"
$ws.Cells.Item(575, 10).Value = 3000
$ws.Cells.Item(575, 11).Value = 600
$ws.Cells.Item(575, 12).Value = 600
$ws.Cells.Item(575, 13).Value = 600
$ws.Cells.Item(575, 16).Value = 120

$ws.Cells.Item(576, 4).Value = 44323
$ws.Cells.Item(576, 9).Value = "Primera"
$ws.Cells.Item(576, 10).Value = 3000
$ws.Cells.Item(576, 11).Value = 500
$ws.Cells.Item(576, 12).Value = 500
$ws.Cells.Item(576, 13).Value = 500
$ws.Cells.Item(576, 16).Value = 100

$ws.Cells.Item(577, 4).Value = 44221
$ws.Cells.Item(577, 9).Value = "Primera"
$ws.Cells.Item(577, 10).Value = 3000
$ws.Cells.Item(577, 11).Value = 550
$ws.Cells.Item(577, 12).Value = 550
$ws.Cells.Item(577, 13).Value = 550
$ws.Cells.Item(577, 16).Value = 110

$ws.Cells.Item(578, 4).Value = 44664
$ws.Cells.Item(578, 9).Value = "Primera"
$ws.Cells.Item(578, 10).Value = 6000
$ws.Cells.Item(578, 11).Value = 600
$ws.Cells.Item(578, 12).Value = 600
$ws.Cells.Item(578, 13).Value = 600
$ws.Cells.Item(578, 16).Value = 120

$ws.Cells.Item(579, 4).Value = 44511
$ws.Cells.Item(579, 9).Value = "Primera"
$ws.Cells.Item(579, 10).Value = 5000
$ws.Cells.Item(579, 11).Value = 500
$ws.Cells.Item(579, 12).Value = 500
$ws.Cells.Item(579, 13).Value = 500
$ws.Cells.Item(579, 16).Value = 100

$ws.Cells.Item(580, 4).Value = 44468
$ws.Cells.Item(580, 9).Value = "Primera"
$ws.Cells.Item(580, 10).Value = 3000
$ws.Cells.Item(580, 11).Value = 650
$ws.Cells.Item(580, 12).Value = 650
$ws.Cells.Item(580, 13).Value = 650
$ws.Cells.Item(580, 16).Value = 130

$ws.Cells.Item(581, 4).Value = 44253
$ws.Cells.Item(581, 9).Value = "Primera"
$ws.Cells.Item(581, 10).Value = 2000
$ws.Cells.Item(581, 11).Value = 600
$ws.Cells.Item(581, 12).Value = 600
$ws.Cells.Item(581, 13).Value = 600
$ws.Cells.Item(581, 16).Value = 120

$ws.Cells.Item(582, 4).Value = 45106
$ws.Cells.Item(582, 9).Value = "Primera"
$ws.Cells.Item(582, 10).Value = 3000
$ws.Cells.Item(582, 11).Value = 650
$ws.Cells.Item(582, 12).Value = 650
$ws.Cells.Item(582, 13).Value = 650
$ws.Cells.Item(582, 16).Value = 130

$ws.Cells.Item(583, 4).Value = 45106
$ws.Cells.Item(583, 9).Value = "Segunda"
$ws.Cells.Item(583, 10).Value = 2000
$ws.Cells.Item(583, 11).Value = 550
$ws.Cells.Item(583, 12).Value = 550
$ws.Cells.Item(583, 13).Value = 550
$ws.Cells.Item(583, 16).Value = 110

$ws.Cells.Item(584, 4).Value = 44169
$ws.Cells.Item(584, 9).Value = "Primera"
$ws.Cells.Item(584, 10).Value = 3000
$ws.Cells.Item(584, 11).Value = 600
$ws.Cells.Item(584, 12).Value = 600
$ws.Cells.Item(584, 13).Value = 600
$ws.Cells.Item(584, 16).Value = 120

$ws.Cells.Item(585, 4).Value = 44358
$ws.Cells.Item(585, 9).Value = "Primera"
$ws.Cells.Item(585, 10).Value = 5000
$ws.Cells.Item(585, 11).Value = 500
$ws.Cells.Item(585, 12).Value = 500
$ws.Cells.Item(585, 13).Value = 500
$ws.Cells.Item(585, 16).Value = 100

$ws.Cells.Item(586, 4).Value = 44235
$ws.Cells.Item(586, 9).Value = "Primera"
$ws.Cells.Item(586, 10).Value = 3000
$ws.Cells.Item(586, 11).Value = 600
$ws.Cells.Item(586, 12).Value = 600
$ws.Cells.Item(586, 13).Value = 600
$ws.Cells.Item(586, 16).Value = 120

$ws.Cells.Item(587, 4).Value = 44391
$ws.Cells.Item(587, 9).Value = "Primera"
$ws.Cells.Item(587, 10).Value = 5000
$ws.Cells.Item(587, 11).Value = 650
$ws.Cells.Item(587, 12).Value = 650
$ws.Cells.Item(587, 13).Value = 650
$ws.Cells.Item(587, 16).Value = 130

$ws.Cells.Item(588, 4).Value = 44634
$ws.Cells.Item(588, 9).Value = "Primera"
$ws.Cells.Item(588, 10).Value = 2000
$ws.Cells.Item(588, 11).Value = 700
$ws.Cells.Item(588, 12).Value = 700
$ws.Cells.Item(588, 13).Value = 700
$ws.Cells.Item(588, 16).Value = 140

$ws.Cells.Item(589, 4).Value = 44634
$ws.Cells.Item(589, 9).Value = "Segunda"
$ws.Cells.Item(589, 10).Value = 2000
$ws.Cells.Item(589, 11).Value = 600
$ws.Cells.Item(589, 12).Value = 600
$ws.Cells.Item(589, 13).Value = 600
$ws.Cells.Item(589, 16).Value = 120

$ws.Cells.Item(590, 4).Value = 44420
$ws.Cells.Item(590, 9).Value = "Primera"
$ws.Cells.Item(590, 10).Value = 3000
$ws.Cells.Item(590, 11).Value = 600
$ws.Cells.Item(590, 12).Value = 600
$ws.Cells.Item(590, 13).Value = 600
$ws.Cells.Item(590, 16).Value = 120

$ws.Cells.Item(591, 4).Value = 44924
$ws.Cells.Item(591, 9).Value = "Primera"
$ws.Cells.Item(591, 10).Value = 3500
$ws.Cells.Item(591, 11).Value = 700
$ws.Cells.Item(591, 12).Value = 700
$ws.Cells.Item(591, 13).Value = 700
$ws.Cells.Item(591, 16).Value = 140

$ws.Cells.Item(592, 4).Value = 44638
$ws.Cells.Item(592, 9).Value = "Primera"
$ws.Cells.Item(592, 10).Value = 4000
$ws.Cells.Item(592, 11).Value = 700
$ws.Cells.Item(592, 12).Value = 700
$ws.Cells.Item(592, 13).Value = 700
$ws.Cells.Item(592, 16).Value = 140

$ws.Cells.Item(593, 4).Value = 45075
$ws.Cells.Item(593, 9).Value = "Primera"
$ws.Cells.Item(593, 10).Value = 3000
$ws.Cells.Item(593, 11).Value = 600
$ws.Cells.Item(593, 12).Value = 600
$ws.Cells.Item(593, 13).Value = 600
$ws.Cells.Item(593, 16).Value = 120

$ws.Cells.Item(594, 4).Value = 45075
$ws.Cells.Item(594, 9).Value = "Segunda"
$ws.Cells.Item(594, 10).Value = 3000
$ws.Cells.Item(594, 11).Value = 500
$ws.Cells.Item(594, 12).Value = 500
$ws.Cells.Item(594, 13).Value = 500
$ws.Cells.Item(594, 16).Value = 100

$ws.Cells.Item(595, 4).Value = 44971
$ws.Cells.Item(595, 9).Value = "Primera"
$ws.Cells.Item(595, 10).Value = 5000
$ws.Cells.Item(595, 11).Value = 550
$ws.Cells.Item(595, 12).Value = 550
$ws.Cells.Item(595, 13).Value = 550
$ws.Cells.Item(595, 16).Value = 110

$ws.Cells.Item(596, 4).Value = 44364
$ws.Cells.Item(596, 9).Value = "Primera"
$ws.Cells.Item(596, 10).Value = 5000
$ws.Cells.Item(596, 11).Value = 500
$ws.Cells.Item(596, 12).Value = 500
$ws.Cells.Item(596, 13).Value = 500
$ws.Cells.Item(596, 16).Value = 100

$ws.Cells.Item(597, 4).Value = 44517
$ws.Cells.Item(597, 9).Value = "Primera"
$ws.Cells.Item(597, 10).Value = 5000
$ws.Cells.Item(597, 11).Value = 550
$ws.Cells.Item(597, 12).Value = 550
$ws.Cells.Item(597, 13).Value = 550
$ws.Cells.Item(597, 16).Value = 110

$ws.Cells.Item(598, 4).Value = 44985
$ws.Cells.Item(598, 9).Value = "Primera"
$ws.Cells.Item(598, 10).Value = 6000
$ws.Cells.Item(598, 11).Value = 500
$ws.Cells.Item(598, 12).Value = 550
$ws.Cells.Item(598, 13).Value = 525
$ws.Cells.Item(598, 16).Value = 105

$ws.Cells.Item(599, 4).Value = 44644
$ws.Cells.Item(599, 9).Value = "Primera"
$ws.Cells.Item(599, 10).Value = 5000
$ws.Cells.Item(599, 11).Value = 700
$ws.Cells.Item(599, 12).Value = 700
$ws.Cells.Item(599, 13).Value = 700
$ws.Cells.Item(599, 16).Value = 140

$ws.Cells.Item(600, 4).Value = 44729
$ws.Cells.Item(600, 9).Value = "Primera"
$ws.Cells.Item(600, 10).Value = 4000
$ws.Cells.Item(600, 11).Value = 700
$ws.Cells.Item(600, 12).Value = 700
$ws.Cells.Item(600, 13).Value = 700
$ws.Cells.Item(600, 16).Value = 140

$ws.Cells.Item(601, 4).Value = 44630
$ws.Cells.Item(601, 9).Value = "Primera"
$ws.Cells.Item(601, 10).Value = 3000
$ws.Cells.Item(601, 11).Value = 800
$ws.Cells.Item(601, 12).Value = 800
$ws.Cells.Item(601, 13).Value = 800
$ws.Cells.Item(601, 16).Value = 160

$ws.Cells.Item(602, 4).Value = 44811
$ws.Cells.Item(602, 9).Value = "Primera"
$ws.Cells.Item(602, 10).Value = 2000
$ws.Cells.Item(602, 11).Value = 1000
$ws.Cells.Item(602, 12).Value = 1000
$ws.Cells.Item(602, 13).Value = 1000
$ws.Cells.Item(602, 16).Value = 200

$ws.Cells.Item(603, 4).Value = 44811
$ws.Cells.Item(603, 9).Value = "Segunda"
$ws.Cells.Item(603, 10).Value = 2000
$ws.Cells.Item(603, 11).Value = 800
$ws.Cells.Item(603, 12).Value = 800
$ws.Cells.Item(603, 13).Value = 800
$ws.Cells.Item(603, 16).Value = 160

$ws.Cells.Item(604, 4).Value = 44167
$ws.Cells.Item(604, 9).Value = "Primera"
$ws.Cells.Item(604, 10).Value = 3000
$ws.Cells.Item(604, 11).Value = 500
$ws.Cells.Item(604, 12).Value = 500
$ws.Cells.Item(604, 13).Value = 500
$ws.Cells.Item(604, 16).Value = 100

$ws.Cells.Item(605, 4).Value = 44637
$ws.Cells.Item(605, 9).Value = "Primera"
$ws.Cells.Item(605, 10).Value = 4000
$ws.Cells.Item(605, 11).Value = 700
$ws.Cells.Item(605, 12).Value = 700
$ws.Cells.Item(605, 13).Value = 700
$ws.Cells.Item(605, 16).Value = 140

$ws.Cells.Item(606, 4).Value = 45112
$ws.Cells.Item(606, 9).Value = "Primera"
$ws.Cells.Item(606, 10).Value = 3000
$ws.Cells.Item(606, 11).Value = 600
$ws.Cells.Item(606, 12).Value = 600
$ws.Cells.Item(606, 13).Value = 600
$ws.Cells.Item(606, 16).Value = 120

$ws.Cells.Item(607, 4).Value = 45112
$ws.Cells.Item(607, 9).Value = "Segunda"
$ws.Cells.Item(607, 10).Value = 2000
$ws.Cells.Item(607, 11).Value = 500
$ws.Cells.Item(607, 12).Value = 500
$ws.Cells.Item(607, 13).Value = 500
$ws.Cells.Item(607, 16).Value = 100

$ws.Cells.Item(608, 4).Value = 45112
$ws.Cells.Item(608, 9).Value = "Primera"
$ws.Cells.Item(608, 10).Value = 5000
$ws.Cells.Item(608, 11).Value = 600
$ws.Cells.Item(608, 12).Value = 600
$ws.Cells.Item(608, 13).Value = 600
$ws.Cells.Item(608, 16).Value = 120

$ws.Cells.Item(609, 4).Value = 44481
$ws.Cells.Item(609, 9).Value = "Primera"
$ws.Cells.Item(609, 10).Value = 4000
$ws.Cells.Item(609, 11).Value = 600
$ws.Cells.Item(609, 12).Value = 600
$ws.Cells.Item(609, 13).Value = 600
$ws.Cells.Item(609, 16).Value = 120

$ws.Cells.Item(610, 4).Value = 44802
$ws.Cells.Item(610, 9).Value = "Primera"
$ws.Cells.Item(610, 10).Value = 3000
$ws.Cells.Item(610, 11).Value = 850
$ws.Cells.Item(610, 12).Value = 850
$ws.Cells.Item(610, 13).Value = 850
$ws.Cells.Item(610, 16).Value = 170

$ws.Cells.Item(611, 4).Value = 44348
$ws.Cells.Item(611, 9).Value = "Primera"
$ws.Cells.Item(611, 10).Value = 6000
$ws.Cells.Item(611, 11).Value = 500
$ws.Cells.Item(611, 12).Value = 500
$ws.Cells.Item(611, 13).Value = 500
$ws.Cells.Item(611, 16).Value = 100

$ws.Cells.Item(612, 4).Value = 44574
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 4000
$ws.Cells.Item(612, 11).Value = 500
$ws.Cells.Item(612, 12).Value = 500
$ws.Cells.Item(612, 13).Value = 500
$ws.Cells.Item(612, 16).Value = 100

$ws.Cells.Item(613, 4).Value = 44566
$ws.Cells.Item(613, 9).Value = "Primera"
$ws.Cells.Item(613, 10).Value = 4000
$ws.Cells.Item(613, 11).Value = 500
$ws.Cells.Item(613, 12).Value = 500
$ws.Cells.Item(613, 13).Value = 500
$ws.Cells.Item(613, 16).Value = 100

$ws.Cells.Item(614, 4).Value = 44594
$ws.Cells.Item(614, 9).Value = "Primera"
$ws.Cells.Item(614, 10).Value = 4000
$ws.Cells.Item(614, 11).Value = 700
$ws.Cells.Item(614, 12).Value = 700
$ws.Cells.Item(614, 13).Value = 700
$ws.Cells.Item(614, 16).Value = 140

$ws.Cells.Item(615, 4).Value = 44749
$ws.Cells.Item(615, 9).Value = "Primera"
$ws.Cells.Item(615, 10).Value = 4000
$ws.Cells.Item(615, 11).Value = 650
$ws.Cells.Item(615, 12).Value = 650
$ws.Cells.Item(615, 13).Value = 650
$ws.Cells.Item(615, 16).Value = 130

$ws.Cells.Item(616, 4).Value = 44957
$ws.Cells.Item(616, 9).Value = "Primera"
$ws.Cells.Item(616, 10).Value = 5000
$ws.Cells.Item(616, 11).Value = 700
$ws.Cells.Item(616, 12).Value = 700
$ws.Cells.Item(616, 13).Value = 700
$ws.Cells.Item(616, 16).Value = 140

$ws.Cells.Item(617, 4).Value = 45089
$ws.Cells.Item(617, 9).Value = "Primera"
$ws.Cells.Item(617, 10).Value = 3000
$ws.Cells.Item(617, 11).Value = 650
$ws.Cells.Item(617, 12).Value = 650
$ws.Cells.Item(617, 13).Value = 650
$ws.Cells.Item(617, 16).Value = 130

$ws.Cells.Item(618, 4).Value = 45089
$ws.Cells.Item(618, 9).Value = "Segunda"
$ws.Cells.Item(618, 10).Value = 2000
$ws.Cells.Item(618, 11).Value = 550
$ws.Cells.Item(618, 12).Value = 550
$ws.Cells.Item(618, 13).Value = 550
$ws.Cells.Item(618, 16).Value = 110

$ws.Cells.Item(619, 4).Value = 44763
$ws.Cells.Item(619, 9).Value = "Primera"
$ws.Cells.Item(619, 10).Value = 5000
$ws.Cells.Item(619, 11).Value = 800
$ws.Cells.Item(619, 12).Value = 800
$ws.Cells.Item(619, 13).Value = 800
$ws.Cells.Item(619, 16).Value = 160

$ws.Cells.Item(620, 4).Value = 44741
$ws.Cells.Item(620, 9).Value = "Primera"
$ws.Cells.Item(620, 10).Value = 5000
$ws.Cells.Item(620, 11).Value = 700
$ws.Cells.Item(620, 12).Value = 700
$ws.Cells.Item(620, 13).Value = 700
$ws.Cells.Item(620, 16).Value = 140

$ws.Cells.Item(621, 4).Value = 45121
$ws.Cells.Item(621, 9).Value = "Primera"
$ws.Cells.Item(621, 10).Value = 5000
$ws.Cells.Item(621, 11).Value = 650
$ws.Cells.Item(621, 12).Value = 650
$ws.Cells.Item(621, 13).Value = 650
$ws.Cells.Item(621, 16).Value = 130

$ws.Cells.Item(622, 4).Value = 44341
$ws.Cells.Item(622, 9).Value = "Primera"
$ws.Cells.Item(622, 10).Value = 4000
$ws.Cells.Item(622, 11).Value = 500
$ws.Cells.Item(622, 12).Value = 500
$ws.Cells.Item(622, 13).Value = 500
$ws.Cells.Item(622, 16).Value = 100

$ws.Cells.Item(623, 4).Value = 44777
$ws.Cells.Item(623, 9).Value = "Primera"
$ws.Cells.Item(623, 10).Value = 3000
$ws.Cells.Item(623, 11).Value = 750
$ws.Cells.Item(623, 12).Value = 750
$ws.Cells.Item(623, 13).Value = 750
$ws.Cells.Item(623, 16).Value = 150

$ws.Cells.Item(624, 4).Value = 44662
$ws.Cells.Item(624, 9).Value = "Primera"
$ws.Cells.Item(624, 10).Value = 5000
$ws.Cells.Item(624, 11).Value = 600
$ws.Cells.Item(624, 12).Value = 600
$ws.Cells.Item(624, 13).Value = 600
$ws.Cells.Item(624, 16).Value = 120

$ws.Cells.Item(625, 4).Value = 44607
$ws.Cells.Item(625, 9).Value = "Primera"
$ws.Cells.Item(625, 10).Value = 3000
$ws.Cells.Item(625, 11).Value = 800
$ws.Cells.Item(625, 12).Value = 800
$ws.Cells.Item(625, 13).Value = 800
$ws.Cells.Item(625, 16).Value = 160

$ws.Cells.Item(626, 4).Value = 45072
$ws.Cells.Item(626, 9).Value = "Primera"
$ws.Cells.Item(626, 10).Value = 3000
$ws.Cells.Item(626, 11).Value = 600
$ws.Cells.Item(626, 12).Value = 600
$ws.Cells.Item(626, 13).Value = 600
$ws.Cells.Item(626, 16).Value = 120

# New row 627 (appended)
$ws.Cells.Item(627, 1).Value = 5
$ws.Cells.Item(627, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(627, 3).Value = "Maule"
$ws.Cells.Item(627, 4).Value = 45072
$ws.Cells.Item(627, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(627, 5).Value = 7
$ws.Cells.Item(627, 6).Value = 100114014
$ws.Cells.Item(627, 7).Value = "Betarraga"
$ws.Cells.Item(627, 8).Value = "Sin especificar"
$ws.Cells.Item(627, 9).Value = "Segunda"
$ws.Cells.Item(627, 10).Value = 2000
$ws.Cells.Item(627, 11).Value = 500
$ws.Cells.Item(627, 12).Value = 500
$ws.Cells.Item(627, 13).Value = 500
$ws.Cells.Item(627, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(627, 15).Value = "Región del Maule"
$ws.Cells.Item(627, 16).Value = 100
$ws.Cells.Item(627, 17).Value = 5
$ws.Cells.Item(627, 18).Value = "Hortaliza"
